# change domain remove whitespace and umlaute
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$domains = @{
    2  = "www.wasim-web.bbq"
    3  = "www.amaanilinda-web.bbq"
    4  = "www.serhii-web.bbq"
    5  = "www.elina-web.bbq"
    6  = "www.mohammed-web.bbq"
    7  = "www.marina-web.bbq"
    8  = "www.ilyas-web.bbq"
    9  = "www.david-web.bbq"
    10 = "www.thomas-web.bbq"
    11 = "www.hadi-web.bbq"
    12 = "www.enes-web.bbq"
    13 = "www.jitender-web.bbq"
    14 = "www.cetin-web.bbq"
    15 = "www.ehsan-web.bbq"
    16 = "www.melzer-web.bbq"
    17 = "www.bogenberg-web.bbq"
}

foreach ($row in $domains.Keys) {
    $ws.Cells.Item($row, 3).Value = $domains[$row]
}
